# Applies cryptocurrency price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.873.06"
$ws.Range("E2").Value = "  +6.36%  "
$ws.Range("D3").Value = "3.009.04"
$ws.Range("E3").Value = "  +3.61%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'582.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "'162.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.79%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.54%  "
$ws.Range("D9").Value = "3.004.50"
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").Value = "'6.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").Value = "'0.156"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.26%  "
$ws.Range("E12").Value = "  +7.30%  "
$ws.Range("E13").Value = "  +9.15%  "
$ws.Range("D14").Value = "'34.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.74%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "65.876.72"
$ws.Range("E16").Value = "  +6.36%  "
$ws.Range("D17").Value = "3.508.19"
$ws.Range("E17").Value = "  +3.45%  "
$ws.Range("E18").Value = "  +7.57%  "
$ws.Range("D19").Value = "3.011.29"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("D20").Value = "'457.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.59%  "
$ws.Range("D21").Value = "'13.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.64%  "
$ws.Range("D22").Value = "'0.689"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.83%  "
$ws.Range("D23").Value = "'7.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.42%  "
$ws.Range("D24").Value = "'82.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.31%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.13%  "
$ws.Range("D26").Value = "'12.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.31%  "
$ws.Range("D27").Value = "'10.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.08%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'8.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.66%  "
$ws.Range("D30").Value = "'2.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +15.85%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "'0.0000103"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.11%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").Value = "'26.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.37%  "
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "'0.993"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.21%  "
$ws.Range("D37").Value = "'5.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.82%  "
$ws.Range("D38").Value = "'2.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.74%  "
$ws.Range("D39").Value = "'49.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'2.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("D41").Value = "'0.309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +15.72%  "
$ws.Range("E42").Value = "  +6.21%  "
$ws.Range("D43").Value = "'43.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.71%  "
$ws.Range("D44").Value = "'8.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("D45").Value = "'390.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.75%  "
$ws.Range("D46").Value = "'0.0355"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.84%  "
$ws.Range("D47").Value = "2.789.53"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "'134.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'23.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.45%  "
$ws.Range("E51").Value = "  +4.04%  "
